# Update "Horarios" workbook: Línea 141 - 303 schedule refresh.
# New scrape timestamp: 02:17:56 (was 01:55:40)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Insert a new data row at row 8, pushing the former rows 8-9 down to 9-10.
$ws1.Rows.Item(8).Insert()

$ws1.Range("A8").Value = "02:17:56"
$ws1.Range("B8").Value = "02:57"
$ws1.Range("C8").Value = "215_ALUAR"
$ws1.Range("D8").Value = 40
$ws1.Range("E8").Value = "LP1912"

# Append a brand-new row at the bottom (row 11).
$ws1.Range("A11").Value = "02:17:56"
$ws1.Range("B11").Value = "04:01"
$ws1.Range("C11").Value = "81_EL PELIGRO"
$ws1.Range("D11").Value = 104
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A2").Value = "Última actualización: 02:17:56"
$ws1.Range("A3").Value = "Total filas: 6"

# ---------------------------------------------------------------
# Sheet 2: "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Insert a new data row at row 7, pushing the former row 7 down to row 8.
$ws2.Rows.Item(7).Insert()

$ws2.Range("A7").Value = "02:17:56"
$ws2.Range("B7").Value = "02:57"
$ws2.Range("C7").Value = "215_ALUAR"
$ws2.Range("D7").Value = 40
$ws2.Range("E7").Value = "LP1912"

$ws2.Range("A2").Value = "Última actualización: 02:17:56"
$ws2.Range("A3").Value = "Total filas: 3"

# ---------------------------------------------------------------
# Sheet 3: "6203-6173"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 02:17:56"
